# DeudoresPrueba.xlsx refresh: re-sorted client list, updated dates/amounts for
# existing debtors, and three newly-added debtor rows (PORTAL ZIPA, SANTANDER
# MADRID, VNZLNO PUNTA ANCA moved down) plus two brand-new clients (CAMPO VERDE
# TOCANCIPA, MERKA FRUVER DEXI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh rows 2-26: client name (B), due date (C), amount owed (D) ---
$ws.Range("B2").Value = "ALISO"
$ws.Range("C2").Value = 45996
$ws.Range("D2").Value = 196000
$ws.Range("B3").Value = "CAMPO VERDE TOCANCIPA"
$ws.Range("C3").Value = 45995
$ws.Range("D3").Value = 635000
$ws.Range("B4").Value = "CAMPO VERDE ZIPAQUIRA"
$ws.Range("C4").Value = 45995
$ws.Range("D4").Value = 684200
$ws.Range("B5").Value = "CARNES JOHANA"
$ws.Range("C5").Value = 45993
$ws.Range("D5").Value = 176000
$ws.Range("B6").Value = "CARNES JOHANA"
$ws.Range("C6").Value = 45996
$ws.Range("D6").Value = 82000
$ws.Range("B7").Value = "CARNILANDIA"
$ws.Range("C7").Value = 45994
$ws.Range("D7").Value = 436700
$ws.Range("B8").Value = "CARNIVOROS"
$ws.Range("C8").Value = 45959
$ws.Range("D8").Value = 200000
$ws.Range("B9").Value = "CIMARRON DORADO"
$ws.Range("C9").Value = 45992
$ws.Range("D9").Value = 407000
$ws.Range("B10").Value = "CIMARRON DORADO"
$ws.Range("C10").Value = 45996
$ws.Range("D10").Value = 298700
$ws.Range("B11").Value = "COCINA CHINA"
$ws.Range("C11").Value = 45994
$ws.Range("D11").Value = 85000
$ws.Range("B12").Value = "DARWIN FUTBOL"
$ws.Range("C12").Value = 45921
$ws.Range("D12").Value = 200000
$ws.Range("B13").Value = "DAVIDCITO"
$ws.Range("C13").Value = 45947
$ws.Range("D13").Value = 100000
$ws.Range("B14").Value = "EL RUBY"
$ws.Range("C14").Value = 45992
$ws.Range("D14").Value = 85100
$ws.Range("B15").Value = "LA PAMPA"
$ws.Range("C15").Value = 45994
$ws.Range("D15").Value = 249000
$ws.Range("B16").Value = "LA SELECTA"
$ws.Range("C16").Value = 45912
$ws.Range("D16").Value = 82000
$ws.Range("B17").Value = "MARIANA"
$ws.Range("C17").Value = 45650
$ws.Range("D17").Value = 171900
$ws.Range("B18").Value = "MERKA FRUVER ALEJANDRO"
$ws.Range("C18").Value = 45995
$ws.Range("D18").Value = 954600
$ws.Range("B19").Value = "MERKA FRUVER DEXI"
$ws.Range("C19").Value = 45995
$ws.Range("D19").Value = 454400
$ws.Range("B20").Value = "NOVILLON SAN MATEO"
$ws.Range("C20").Value = 45971
$ws.Range("D20").Value = 83000
$ws.Range("B21").Value = "PINILLA"
$ws.Range("C21").Value = 45931
$ws.Range("D21").Value = 166000
$ws.Range("B22").Value = "PINILLA"
$ws.Range("C22").Value = 45924
$ws.Range("D22").Value = 16000
$ws.Range("B23").Value = "PINILLA SOACHA"
$ws.Range("C23").Value = 45993
$ws.Range("D23").Value = 129000
$ws.Range("B24").Value = "PLAZA JESSICA"
$ws.Range("C24").Value = 45993
$ws.Range("D24").Value = 621000
$ws.Range("B25").Value = "PLAZA JESSICA"
$ws.Range("C25").Value = 45995
$ws.Range("D25").Value = 1580300
$ws.Range("B26").Value = "PORTAL ZIPA"
$ws.Range("C26").Value = 45995
$ws.Range("D26").Value = 66400

# --- Append newly-added debtor rows 27-29 ---
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "SANTANDER MADRID"
$ws.Range("C27").Value = 45996
$ws.Range("D27").Value = 63000
$ws.Range("E27").Value = $false
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "SANTANDER SUR"
$ws.Range("C28").Value = 45993
$ws.Range("D28").Value = 80000
$ws.Range("E28").Value = $false
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "VNZLNO PUNTA ANCA"
$ws.Range("C29").Value = 45992
$ws.Range("D29").Value = 82000
$ws.Range("E29").Value = $false

# New date cells should carry the same custom date format already used by column C
$ws.Range("C27:C29").NumberFormat = "YYYY-MM-DD"

# Widen the Cliente/Fecha columns to fit the refreshed content (matches the
# "best fit" column widths Excel applies automatically after a data refresh)
$ws.Columns("B:B").AutoFit()
$ws.Columns("C:C").AutoFit()

# Leave the selection where the author left it after entering the new rows
$null = $ws.Range("K27").Select()
